$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update David's (row 7) Chinese score (C7) from 90 to 82
$ws.Range("C7").Value = 82

# Replace the literal totals in I7/J7 with formulas
$ws.Range("I7").Formula = "=SUM(C7:H7)"
$ws.Range("J7").Formula = "=AVERAGE(C7:H7)"
